$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "runs" (C) and "balls" (D) columns were updated with new activity figures.
# Values are kept as text (matching the sheet's existing numberStoredAsText
# convention) by prefixing with an apostrophe, Excel's standard "force text"
# entry marker.
$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "'3"

$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'0"

$ws.Range("C4").Value = "'0"
$ws.Range("D4").Value = "'0"

$ws.Range("D5").Value = "'2"

$ws.Range("C6").Value = "'3"
$ws.Range("D6").Value = "'1"
